# Update the title from a declarative statement to a question, and move the
# "_GoBack" bookmark (Word's "last edit position" marker) from its old spot
# near the equations to the new edit point inside the title.

$d = $word.ActiveDocument

# Title paragraph is paragraph 1: "Matrix neural network is the alternative
# of convolutional neural network?" The first run holds everything up to
# (but excluding) the trailing "?", which lives in its own run.
$titleRun = $d.Paragraphs.Item(1).Range
$oldText = "Matrix neural network is the alternative of convolutional neural network"
$boundary = $oldText.Length

# Guard the boundary with the old "?" run so the upcoming text replacement
# does not get merged into it.
$guard = $d.Bookmarks.Add("ZZZGuard", $d.Range($boundary, $boundary))

# Replace the sentence with the new wording (same length: 72 characters).
$newText = "Is matrix neural network the alternative of convolutional neural network"
$d.Range(0, $boundary).Text = $newText

# Force a run split right after "Is " (position 3) with a scratch bookmark,
# then drop "_GoBack" right after "Is m" (position 4) -- adding a bookmark
# named "_GoBack" automatically relocates it from wherever it used to be
# (next to the equations further down) to this new spot.
$splitPoint = $d.Bookmarks.Add("ZZZSplit", $d.Range(3, 3))
$goBack = $d.Bookmarks.Add("_GoBack", $d.Range(4, 4))

# Clean up the scratch bookmarks (keeps the run split, removes the markers).
$d.Bookmarks.Item("ZZZSplit").Delete()
$d.Bookmarks.Item("ZZZGuard").Delete()
